$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: nudge the window position to match the author's session
# (xWindow stays 240, yWindow moves from 11505 to 12105). Window chrome
# state may not round-trip through every host, but set it anyway.
try {
    $win = $excel.ActiveWindow
    $win.Top = 12105
    $win.Left = 240
} catch {
}

# The "MLC" tag used for column F (USE_ACTUAL_MODEL) was renamed to the
# more specific "PR_C_Y2" across all data rows.
$ws.Range("F2").Value = "PR_C_Y2"
$ws.Range("F3").Value = "PR_C_Y2"
$ws.Range("F4").Value = "PR_C_Y2"

# Move/update the active selection from F3 to H3.
$ws.Range("H3").Select()
